$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing TEXT storage.
# Several "Price" values in column D are plain decimals (e.g. "112.01"),
# and Excel auto-converts a bare numeric-looking string assigned via
# .Value into a real number. The source data must stay text (as in the
# original workbook), so we flip the cell to the "@" (Text) number
# format just long enough to assign the string, then restore the
# cell's original Style so no stray formatting is left behind.
function Set-TextValue($rng, [string]$val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "43.357.62"
$ws.Range("E2").Value = "  -0.84%  "

# Row 3
Set-TextValue $ws.Range("D3") "2.277.74"
$ws.Range("E3").Value = "  -0.47%  "

# Row 4
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
Set-TextValue $ws.Range("D5") "112.01"
$ws.Range("E5").Value = "  +1.09%  "

# Row 6
Set-TextValue $ws.Range("D6") "264.55"
$ws.Range("E6").Value = "  -1.02%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.644"
$ws.Range("E7").Value = "  +3.11%  "

# Row 8
$ws.Range("E8").Value = "  +0.02%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.609"
$ws.Range("E9").Value = "  -1.69%  "

# Row 10
Set-TextValue $ws.Range("D10") "46.58"
$ws.Range("E10").Value = "  -2.04%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.0934"
$ws.Range("E11").Value = "  -1.08%  "

# Row 12
Set-TextValue $ws.Range("D12") "9.18"
$ws.Range("E12").Value = "  +3.78%  "

# Row 13
$ws.Range("E13").Value = "  +1.51%  "

# Row 14
Set-TextValue $ws.Range("D14") "15.34"
$ws.Range("E14").Value = "  -2.57%  "

# Row 15
Set-TextValue $ws.Range("D15") "2.619.09"
$ws.Range("E15").Value = "  -0.50%  "

# Row 16
$ws.Range("E16").Value = "  +1.94%  "

# Row 17
Set-TextValue $ws.Range("D17") "2.265.56"
$ws.Range("E17").Value = "  -1.18%  "

# Row 18
Set-TextValue $ws.Range("D18") "43.211.79"
$ws.Range("E18").Value = "  -0.88%  "

# Row 19
$ws.Range("E19").Value = "  -1.65%  "

# Row 20
Set-TextValue $ws.Range("D20") "6.76"
$ws.Range("E20").Value = "  +2.40%  "

# Row 21
Set-TextValue $ws.Range("D21") "72.17"
$ws.Range("E21").Value = "  -0.27%  "

# Row 22
$ws.Range("E22").Value = "  -1.03%  "

# Row 23
Set-TextValue $ws.Range("D23") "234.46"
$ws.Range("E23").Value = "  +0.92%  "

# Row 24
Set-TextValue $ws.Range("D24") "2.88"
$ws.Range("E24").Value = "  +2.97%  "

# Row 25
Set-TextValue $ws.Range("D25") "9.34"
$ws.Range("E25").Value = "  -2.91%  "

# Row 26
$ws.Range("E26").Value = "  +2.00%  "

# Row 27
Set-TextValue $ws.Range("D27") "11.34"
$ws.Range("E27").Value = "  -2.28%  "

# Row 28
Set-TextValue $ws.Range("D28") "41.43"
$ws.Range("E28").Value = "  -0.82%  "

# Row 29
Set-TextValue $ws.Range("D29") "3.35"
$ws.Range("E29").Value = "  -1.66%  "

# Row 30
$ws.Range("E30").Value = "  -0.87%  "

# Row 31
Set-TextValue $ws.Range("D31") "173.47"
$ws.Range("E31").Value = "  -1.51%  "

# Row 32
Set-TextValue $ws.Range("D32") "21.49"
$ws.Range("E32").Value = "  -0.25%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.0897"
$ws.Range("E33").Value = "  -3.44%  "

# Row 34
Set-TextValue $ws.Range("D34") "5.66"
$ws.Range("E34").Value = "  +0.50%  "

# Row 35
$ws.Range("E35").Value = "  +3.12%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.0378"
$ws.Range("E36").Value = "  +4.83%  "

# Row 37
Set-TextValue $ws.Range("D37") "4.67"
$ws.Range("E37").Value = "  -0.73%  "

# Row 38
Set-TextValue $ws.Range("D38") "3.88"
$ws.Range("E38").Value = "  +3.06%  "

# Row 39
$ws.Range("E39").Value = "  -2.92%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.59"
$ws.Range("E40").Value = "  +7.99%  "

# Row 41
Set-TextValue $ws.Range("D41") "14.27"
$ws.Range("E41").Value = "  +3.89%  "

# Row 42
Set-TextValue $ws.Range("D42") "75.54"
$ws.Range("E42").Value = "  +5.99%  "

# Row 43
$ws.Range("E43").Value = "  -2.78%  "

# Row 44
Set-TextValue $ws.Range("D44") "6.10"
$ws.Range("E44").Value = "  -1.24%  "

# Row 45
$ws.Range("E45").Value = "  -0.08%  "

# Row 46
$ws.Range("E46").Value = "  -2.78%  "

# Row 47
Set-TextValue $ws.Range("D47") "1.27"
$ws.Range("E47").Value = "  +3.96%  "

# Row 48
Set-TextValue $ws.Range("D48") "8.56"
$ws.Range("E48").Value = "  -3.52%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.0992"
$ws.Range("E49").Value = "  -1.47%  "

# Row 50
Set-TextValue $ws.Range("D50") "100.48"
$ws.Range("E50").Value = "  -0.98%  "

# Row 51
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-TextValue $ws.Range("D51") "0.434"
$ws.Range("E51").Value = "  -2.92%  "
